$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1128.2
$ws.Cells.Item(15, 9).Value = 1128.2
$ws.Cells.Item(15, 11).Value = 3384.6
$ws.Cells.Item(15, 13).Value = -3215.6

$ws.Cells.Item(63, 8).Value = 30866.334
$ws.Cells.Item(63, 10).Value = 30866.334
$ws.Cells.Item(63, 12).Value = 30866.334
$ws.Cells.Item(63, 14).Value = -32114.334

$ws.Cells.Item(64, 8).Value = 3117.1428
$ws.Cells.Item(64, 9).Value = 3200
$ws.Cells.Item(64, 10).Value = 3099.1304
$ws.Cells.Item(64, 11).Value = 3200
$ws.Cells.Item(64, 12).Value = 3099.1304
$ws.Cells.Item(64, 13).Value = -2952
$ws.Cells.Item(64, 14).Value = -3595.1304

$ws.Cells.Item(66, 8).Value = 30866.334
$ws.Cells.Item(66, 10).Value = 30866.334
$ws.Cells.Item(66, 12).Value = 92599.00199999999
$ws.Cells.Item(66, 14).Value = -98839.00199999999

$ws.Cells.Item(67, 8).Value = 3117.1428
$ws.Cells.Item(67, 9).Value = 3200
$ws.Cells.Item(67, 10).Value = 3099.1304
$ws.Cells.Item(67, 11).Value = 3200
$ws.Cells.Item(67, 12).Value = 3099.1304
$ws.Cells.Item(67, 13).Value = -2342
$ws.Cells.Item(67, 14).Value = -4815.1304

$ws.Cells.Item(76, 8).Value = 9429.208000000001
$ws.Cells.Item(76, 9).Value = 18975.125
$ws.Cells.Item(76, 10).Value = 4656.25
$ws.Cells.Item(76, 11).Value = 18975.125
$ws.Cells.Item(76, 12).Value = 4656.25
$ws.Cells.Item(76, 13).Value = -18660.125
$ws.Cells.Item(76, 14).Value = -5286.25

$ws.Cells.Item(79, 8).Value = 9429.208000000001
$ws.Cells.Item(79, 9).Value = 18975.125
$ws.Cells.Item(79, 10).Value = 4656.25
$ws.Cells.Item(79, 11).Value = 18975.125
$ws.Cells.Item(79, 12).Value = 4656.25
$ws.Cells.Item(79, 13).Value = -17883.125
$ws.Cells.Item(79, 14).Value = -6840.25

$ws.Cells.Item(107, 8).Value = 1075.45
$ws.Cells.Item(107, 9).Value = 1115.6154
$ws.Cells.Item(107, 11).Value = 1115.6154
$ws.Cells.Item(107, 13).Value = 804.3846000000001

$ws.Cells.Item(111, 8).Value = 3451.8
$ws.Cells.Item(111, 9).Value = 2643.3635
$ws.Cells.Item(111, 11).Value = 7930.0905
$ws.Cells.Item(111, 13).Value = -4863.0905

$ws.Cells.Item(112, 8).Value = 1241
$ws.Cells.Item(112, 9).Value = 778
$ws.Cells.Item(112, 10).Value = 1333.6
$ws.Cells.Item(112, 11).Value = 2334
$ws.Cells.Item(112, 12).Value = 4000.8
$ws.Cells.Item(112, 13).Value = -1226
$ws.Cells.Item(112, 14).Value = -6216.799999999999

$ws.Cells.Item(115, 8).Value = 2178.7778
$ws.Cells.Item(115, 9).Value = 646.44446
$ws.Cells.Item(115, 11).Value = 1939.33338
$ws.Cells.Item(115, 13).Value = -372.33338

$ws.Cells.Item(116, 8).Value = 5086.346
$ws.Cells.Item(116, 9).Value = 4838.9287
$ws.Cells.Item(116, 10).Value = 5375
$ws.Cells.Item(116, 11).Value = 4838.9287
$ws.Cells.Item(116, 12).Value = 5375
$ws.Cells.Item(116, 13).Value = -1396.9287
$ws.Cells.Item(116, 14).Value = -12259

$ws.Cells.Item(138, 8).Value = 3398.98
$ws.Cells.Item(138, 9).Value = 1714.4584
$ws.Cells.Item(138, 10).Value = 3930.9343
$ws.Cells.Item(138, 11).Value = 5143.3752
$ws.Cells.Item(138, 12).Value = 11792.8029
$ws.Cells.Item(138, 13).Value = -3.375200000000405
$ws.Cells.Item(138, 14).Value = -22072.8029

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 28447.361
$ws.Cells.Item(32, 9).Value = 10529.953
$ws.Cells.Item(32, 10).Value = 88800.734
$ws.Cells.Item(32, 11).Value = 10529.953
$ws.Cells.Item(32, 12).Value = 88800.734
$ws.Cells.Item(32, 13).Value = -10242.953
$ws.Cells.Item(32, 14).Value = -89374.734

$ws.Cells.Item(74, 8).Value = 4918.478
$ws.Cells.Item(74, 9).Value = 4666.9287
$ws.Cells.Item(74, 11).Value = 4666.9287
$ws.Cells.Item(74, 13).Value = -3792.9287

$ws.Cells.Item(77, 8).Value = 4918.478
$ws.Cells.Item(77, 9).Value = 4666.9287
$ws.Cells.Item(77, 11).Value = 23334.6435
$ws.Cells.Item(77, 13).Value = -18966.6435

$ws.Cells.Item(118, 8).Value = 38632.8
$ws.Cells.Item(118, 10).Value = 38632.8
$ws.Cells.Item(118, 12).Value = 38632.8
$ws.Cells.Item(118, 14).Value = -41946.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(51, 8).Value = 14933.333
$ws.Cells.Item(51, 10).Value = 14933.333
$ws.Cells.Item(51, 12).Value = 14933.333
$ws.Cells.Item(51, 14).Value = -15915.333

$ws.Cells.Item(107, 8).Value = 4132.952
$ws.Cells.Item(107, 9).Value = 4251.5264
$ws.Cells.Item(107, 11).Value = 4251.5264
$ws.Cells.Item(107, 13).Value = -2331.5264

$ws.Cells.Item(134, 8).Value = 2013
$ws.Cells.Item(134, 9).Value = 1584.8334
$ws.Cells.Item(134, 10).Value = 2869.3333
$ws.Cells.Item(134, 11).Value = 4754.5002
$ws.Cells.Item(134, 12).Value = 8607.999899999999
$ws.Cells.Item(134, 13).Value = -2219.5002
$ws.Cells.Item(134, 14).Value = -13677.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3252.7869
$ws.Cells.Item(31, 9).Value = 1518.2927
$ws.Cells.Item(31, 10).Value = 6808.5
$ws.Cells.Item(31, 11).Value = 1518.2927
$ws.Cells.Item(31, 12).Value = 6808.5
$ws.Cells.Item(31, 13).Value = -1223.2927
$ws.Cells.Item(31, 14).Value = -7398.5

$ws.Cells.Item(34, 8).Value = 3252.7869
$ws.Cells.Item(34, 9).Value = 1518.2927
$ws.Cells.Item(34, 10).Value = 6808.5
$ws.Cells.Item(34, 11).Value = 1518.2927
$ws.Cells.Item(34, 12).Value = 6808.5
$ws.Cells.Item(34, 13).Value = -1316.2927
$ws.Cells.Item(34, 14).Value = -7212.5

$ws.Cells.Item(58, 8).Value = 1678.7368
$ws.Cells.Item(58, 9).Value = 1727.5555
$ws.Cells.Item(58, 10).Value = 800
$ws.Cells.Item(58, 11).Value = 1727.5555
$ws.Cells.Item(58, 12).Value = 800
$ws.Cells.Item(58, 13).Value = -1524.5555
$ws.Cells.Item(58, 14).Value = -1206

$ws.Cells.Item(62, 8).Value = 2833.3333
$ws.Cells.Item(62, 10).Value = 4000
$ws.Cells.Item(62, 12).Value = 4000
$ws.Cells.Item(62, 14).Value = -5248

$ws.Cells.Item(65, 8).Value = 2833.3333
$ws.Cells.Item(65, 10).Value = 4000
$ws.Cells.Item(65, 12).Value = 20000
$ws.Cells.Item(65, 14).Value = -26240

$ws.Cells.Item(134, 8).Value = 4629
$ws.Cells.Item(134, 9).Value = 5394.7036
$ws.Cells.Item(134, 10).Value = 1183.3334
$ws.Cells.Item(134, 11).Value = 16184.1108
$ws.Cells.Item(134, 12).Value = 3550.0002
$ws.Cells.Item(134, 13).Value = -13649.1108
$ws.Cells.Item(134, 14).Value = -8620.0002

$ws.Cells.Item(135, 8).Value = 60520
$ws.Cells.Item(135, 10).Value = 60520
$ws.Cells.Item(135, 12).Value = 60520
$ws.Cells.Item(135, 14).Value = -70660

$ws.Cells.Item(136, 8).Value = 1678.7368
$ws.Cells.Item(136, 9).Value = 1727.5555
$ws.Cells.Item(136, 10).Value = 800
$ws.Cells.Item(136, 11).Value = 5182.666499999999
$ws.Cells.Item(136, 12).Value = 2400
$ws.Cells.Item(136, 13).Value = -2632.666499999999
$ws.Cells.Item(136, 14).Value = -7500

$ws.Cells.Item(138, 8).Value = 41978.4
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 41978.4
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 41978.4
$ws.Cells.Item(138, 14).Value = -52258.4
$ws.Cells.Item(138, 13).ClearContents()

$ws.Cells.Item(140, 8).Value = 78000
$ws.Cells.Item(140, 10).Value = 78000
$ws.Cells.Item(140, 12).Value = 78000
$ws.Cells.Item(140, 14).Value = -88360

$ws.Cells.Item(141, 8).Value = 46899.43
$ws.Cells.Item(141, 10).Value = 51333.332
$ws.Cells.Item(141, 12).Value = 51333.332
$ws.Cells.Item(141, 14).Value = -61693.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(104, 8).Value = 66670268
$ws.Cells.Item(104, 10).Value = 66670268
$ws.Cells.Item(104, 12).Value = 200010804
$ws.Cells.Item(104, 14).Value = -200016046

$ws.Cells.Item(131, 8).Value = 921.0707
$ws.Cells.Item(131, 9).Value = 488.33334
$ws.Cells.Item(131, 10).Value = 948.9892599999999
$ws.Cells.Item(131, 11).Value = 1465.00002
$ws.Cells.Item(131, 12).Value = 2846.96778
$ws.Cells.Item(131, 13).Value = 3574.99998
$ws.Cells.Item(131, 14).Value = -12926.96778

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 34861.91
$ws.Cells.Item(70, 9).Value = 47579.5
$ws.Cells.Item(70, 10).Value = 4339.7
$ws.Cells.Item(70, 11).Value = 47579.5
$ws.Cells.Item(70, 12).Value = 4339.7
$ws.Cells.Item(70, 13).Value = -47309.5
$ws.Cells.Item(70, 14).Value = -4879.7

$ws.Cells.Item(73, 8).Value = 34861.91
$ws.Cells.Item(73, 9).Value = 47579.5
$ws.Cells.Item(73, 10).Value = 4339.7
$ws.Cells.Item(73, 11).Value = 47579.5
$ws.Cells.Item(73, 12).Value = 4339.7
$ws.Cells.Item(73, 13).Value = -46643.5
$ws.Cells.Item(73, 14).Value = -6211.7

$ws.Cells.Item(113, 8).Value = 1580.6364
$ws.Cells.Item(113, 9).Value = 1446.3684
$ws.Cells.Item(113, 10).Value = 2431
$ws.Cells.Item(113, 11).Value = 1446.3684
$ws.Cells.Item(113, 12).Value = 2431
$ws.Cells.Item(113, 13).Value = 723.6315999999999
$ws.Cells.Item(113, 14).Value = -6771

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3053.1177
$ws.Cells.Item(61, 9).Value = 2075.25
$ws.Cells.Item(61, 11).Value = 2075.25
$ws.Cells.Item(61, 13).Value = -1873.25

$ws.Cells.Item(113, 8).Value = 3053.1177
$ws.Cells.Item(113, 9).Value = 2075.25
$ws.Cells.Item(113, 11).Value = 2075.25
$ws.Cells.Item(113, 13).Value = 94.75

$ws.Cells.Item(122, 8).Value = 4693.778
$ws.Cells.Item(122, 9).Value = 4051.3333
$ws.Cells.Item(122, 10).Value = 5015
$ws.Cells.Item(122, 11).Value = 12153.9999
$ws.Cells.Item(122, 12).Value = 15045
$ws.Cells.Item(122, 13).Value = -9703.999899999999
$ws.Cells.Item(122, 14).Value = -19945

$ws.Cells.Item(132, 8).Value = 3095
$ws.Cells.Item(132, 9).Value = 2553.6
$ws.Cells.Item(132, 11).Value = 7660.799999999999
$ws.Cells.Item(132, 13).Value = -5130.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 3564.3
$ws.Cells.Item(132, 9).Value = 4171.8237
$ws.Cells.Item(132, 10).Value = 2769.8462
$ws.Cells.Item(132, 11).Value = 12515.4711
$ws.Cells.Item(132, 12).Value = 8309.5386
$ws.Cells.Item(132, 13).Value = -9985.471099999999
$ws.Cells.Item(132, 14).Value = -13369.5386

$ws.Cells.Item(136, 8).Value = 1564.3636
$ws.Cells.Item(136, 9).Value = 1317.7317
$ws.Cells.Item(136, 10).Value = 4935
$ws.Cells.Item(136, 11).Value = 3953.1951
$ws.Cells.Item(136, 12).Value = 14805
$ws.Cells.Item(136, 13).Value = -1403.1951
$ws.Cells.Item(136, 14).Value = -19905
